# Update cryptos list values (price + 1h volume change) to reflect latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'30.411.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.25%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.102.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.20%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'334.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.38%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +0.05%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.5217"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.08%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.4559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.58%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'54.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +13.68%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.08891"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.05%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +1.02%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'24.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.16%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'2.098.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.33%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.797"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.59%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'8.030"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.32%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'97.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.50%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.00001145"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.03%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.07%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.06626"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.19%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'19.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.60%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  -0.01%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.295"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.37%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'30.487.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.20%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'12.33"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'2.356"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.30%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.333.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.77%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'22.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.44%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'162.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.34%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -4.95%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'133.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.15%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.205"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.93%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.1067"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.55%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.657"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.75%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'6.398"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.53%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'3.931"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.11%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'10.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.59%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'5.870"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.54%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.02573"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.49%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  +1.74%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.2321"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.68%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'12.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.56%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.6869"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.19%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.247"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.72%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'2.322"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.62%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.6393"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.41%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'14.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.24%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'3.662"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.87%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.247"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.63%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'BabyDogeCoin"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.00000000343"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +14.02%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'Aave"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'83.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.80%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'WEMIXTOKEN"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'1.202"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.00%  "
$ws.Range("E51").Style = "Normal"
